$d = $word.ActiveDocument

$d.Content.Find.Execute("I wil now move object A", $true, $false, $false, $false, $false, $true, 1, $false, "I will now move object A", 2)
$d.Content.Find.Execute("Leaving A anC together", $true, $false, $false, $false, $false, $true, 1, $false, "Leaving A and C together", 2)
